# Target doc (6 paragraphs):
#   1: <Tag1><Tag2> <TestTag>
#   2: (empty)
#   3: <TagInNewLine>
#   4: (empty)                              -> becomes "<Tag with spaces>" + _GoBack bookmark
#   5: <Adress>  (spell-check split runs)   -> unchanged
#   6: <Adress>  (single run) + _GoBack     -> bookmark removed, text unchanged
#
# Net effect: the _GoBack bookmark moves from the last paragraph to the
# (formerly empty) 4th paragraph, which also gains the "<Tag with spaces>" text.

$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark from the last paragraph first (while
# it is still the only one in the document), so it doesn't linger once a
# second bookmark of the same name is introduced below.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Turn the empty 4th paragraph into "<Tag with spaces>" and give it the
# _GoBack bookmark, by replacing its contents with the literal OOXML.
$p4 = $d.Paragraphs.Item(4)
$p4xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:r><w:t>&lt;Tag with spaces&gt;</w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
         '<w:bookmarkEnd w:id="0"/>' +
         '</w:p>'
[void]$p4.Range.InsertXML($p4xml)
